$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12970
$ws.Range("C3:C5").Value = 10980
$ws.Range("C6:C14").Value = 9912
$ws.Range("C15:C18").Value = 9457
$ws.Range("C19:C30").Value = 8962
$ws.Range("C31:C41").Value = 8661
$ws.Range("C42:C60").Value = 8445
$ws.Range("C61:C63").Value = 8173
$ws.Range("C64:C94").Value = 8120
$ws.Range("C95:C98").Value = 7855
$ws.Range("C99:C217").Value = 7808
$ws.Range("C218:C252").Value = 7569
